$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.413.74"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.436.78"
$ws.Range("E3").Value = "  -4.71%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "592.44"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "136.50"
$ws.Range("E6").Value = "  -8.19%  "
$ws.Range("D7").Value = "3.437.88"
$ws.Range("E7").Value = "  -4.63%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Value = "7.33"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").Value = "0.122"
$ws.Range("E11").Value = "  -10.09%  "
$ws.Range("D12").Value = "0.382"
$ws.Range("E12").Value = "  -7.69%  "
$ws.Range("D13").Value = "4.027.30"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").Value = "0.0000184"
$ws.Range("E14").Value = "  -11.19%  "
$ws.Range("D15").Value = "26.77"
$ws.Range("E15").Value = "  -10.06%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.115"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.431.68"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").Value = "65.338.13"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").Value = "10.15"
$ws.Range("E19").Value = "  -10.17%  "
$ws.Range("D20").Value = "5.76"
$ws.Range("E20").Value = "  -9.03%  "
$ws.Range("D21").Value = "13.68"
$ws.Range("E21").Value = "  -8.47%  "
$ws.Range("D22").Value = "399.36"
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -10.15%  "
$ws.Range("D24").Value = "73.85"
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.587.84"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").Value = "  -13.79%  "
$ws.Range("D28").Value = "2.34"
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  -12.64%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  -11.52%  "
$ws.Range("D32").Value = "3.450.75"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "0.148"
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("D35").Value = "22.82"
$ws.Range("E35").Value = "  -9.50%  "
$ws.Range("D36").Value = "1.24"
$ws.Range("E36").Value = "  -15.07%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "172.71"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  -10.65%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  -9.38%  "
$ws.Range("D40").Value = "4.84"
$ws.Range("E40").Value = "  -13.87%  "
$ws.Range("D41").Value = "0.0765"
$ws.Range("E41").Value = "  -10.56%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "44.55"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.808"
$ws.Range("E43").Value = "  -9.53%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "4.54"
$ws.Range("E44").Value = "  -12.73%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "1.63"
$ws.Range("E46").Value = "  -12.64%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "6.62"
$ws.Range("E48").Value = "  -7.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "21.60"
$ws.Range("E49").Value = "  -9.01%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -15.27%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.216.80"
$ws.Range("E51").Value = "  -7.87%  "
